$wb = $excel.ActiveWorkbook

# --- StatOutput sheet: update the summary counts row ---
$statOutput = $wb.Worksheets.Item("StatOutput")

# Force these cells to be stored as text (matching the original workbook where
# these numeric-looking counts are stored as shared strings, not numbers).
$statOutput.Range("A2").NumberFormat = "@"
$statOutput.Range("A2").Value = "0"

$statOutput.Range("B2").NumberFormat = "@"
$statOutput.Range("B2").Value = "0"

$statOutput.Range("C2").NumberFormat = "@"
$statOutput.Range("C2").Value = "2"

# D2 (number_of_study) stays "1" - unchanged.

# --- StatOutput_Message sheet: update the Cypher query text (row 18) ---
$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")

$newQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Border Collie']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$statOutputMessage.Range("A18").Value = $newQuery
